# "loop control statements added"
#
# Original paragraph (single run):
#   -- Break and continue
#
# New paragraph (six runs, all sharing the same rPr: sz=32 / szCs=32):
#   "-- " | "Loop Control Statements – " | "Break" | ", " | "continue" | ", pass"

$d = $word.ActiveDocument

# Locate the paragraph robustly via Find (independent of paragraph index).
$hit = $d.Content
$found = $hit.Find.Execute("-- Break and continue", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph text '-- Break and continue'"
}

$paraStart = $hit.Start

# Build the full replacement text in one go (keeps it inside a single run
# with the original formatting intact), then split that run into the six
# pieces shown in the diff.
$newText = "-- Loop Control Statements " + [char]0x2013 + " Break, continue, pass"
$hit.Text = $newText

$paraEnd = $paraStart + $newText.Length

# Character offsets (relative to $paraStart) where a new run must begin.
#   0  -> "-- "
#   3  -> "Loop Control Statements – "
#   29 -> "Break"
#   34 -> ", "
#   36 -> "continue"
#   44 -> ", pass"
$splitOffsets = @(3, 29, 34, 36, 44)

foreach ($offset in $splitOffsets) {
    $splitPos = $paraStart + $offset
    # Selecting from the split point to the (fixed) paragraph end and
    # toggling a character property forces the engine to break this span
    # off into its own run while preserving the inherited run formatting
    # (sz/szCs) on every resulting run.
    $tail = $d.Range($splitPos, $paraEnd)
    $tail.Bold = 1
    $tail.Bold = 0
}
